# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Lamia Profits workbook (updates H/I/J/K/L/M/N columns per sheet/row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 877.7895
$ws.Range("I12").Value = 165.26666
$ws.Range("K12").Value = 165.26666
$ws.Range("M12").Value = 4.733339999999998

$ws.Range("H31").Value = 749
$ws.Range("I31").Value = 998
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 2994
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = -2764
$ws.Range("N31").Value = -1960

$ws.Range("H41").Value = 3089
$ws.Range("I41").Value = 3327.5715
$ws.Range("J41").Value = 2532.3333
$ws.Range("K41").Value = 3327.5715
$ws.Range("L41").Value = 2532.3333
$ws.Range("M41").Value = -2887.5715
$ws.Range("N41").Value = -3412.3333

$ws.Range("H53").Value = 1610.5883
$ws.Range("I53").Value = 356.8
$ws.Range("K53").Value = 356.8
$ws.Range("M53").Value = 280.2

$ws.Range("H62").Value = 3477.6287
$ws.Range("I62").Value = 2505.3794
$ws.Range("K62").Value = 2505.3794
$ws.Range("M62").Value = -1881.3794

$ws.Range("H64").Value = 8333.833000000001
$ws.Range("J64").Value = 10001
$ws.Range("L64").Value = 10001
$ws.Range("N64").Value = -10497

$ws.Range("H65").Value = 3477.6287
$ws.Range("I65").Value = 2505.3794
$ws.Range("K65").Value = 12526.897
$ws.Range("M65").Value = -9406.896999999999

$ws.Range("H67").Value = 8333.833000000001
$ws.Range("J67").Value = 10001
$ws.Range("L67").Value = 10001
$ws.Range("N67").Value = -11717

$ws.Range("H88").Value = 23544.545
$ws.Range("J88").Value = 26554.555
$ws.Range("L88").Value = 26554.555
$ws.Range("N88").Value = -27366.555

$ws.Range("H91").Value = 23544.545
$ws.Range("J91").Value = 26554.555
$ws.Range("L91").Value = 26554.555
$ws.Range("N91").Value = -29362.555

$ws.Range("H112").Value = 2015.5714
$ws.Range("I112").Value = 1246
$ws.Range("K112").Value = 3738
$ws.Range("M112").Value = -2630

$ws.Range("H113").Value = 7284.353
$ws.Range("I113").Value = 4450
$ws.Range("J113").Value = 8156.4614
$ws.Range("K113").Value = 4450
$ws.Range("L113").Value = 8156.4614
$ws.Range("M113").Value = -1196
$ws.Range("N113").Value = -14664.4614

$ws.Range("H138").Value = 3039.818
$ws.Range("J138").Value = 3533.1667
$ws.Range("L138").Value = 10599.5001
$ws.Range("N138").Value = -20879.5001


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 21847.2
$ws.Range("I2").Value = 1897.4
$ws.Range("K2").Value = 1897.4
$ws.Range("M2").Value = -1784.4

$ws.Range("H32").Value = 1828.6812
$ws.Range("I32").Value = 1430.2307
$ws.Range("K32").Value = 1430.2307
$ws.Range("M32").Value = -1143.2307

$ws.Range("H116").Value = 21847.2
$ws.Range("I116").Value = 1897.4
$ws.Range("K116").Value = 1897.4
$ws.Range("M116").Value = 396.5999999999999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 21847.2
$ws.Range("I3").Value = 1897.4
$ws.Range("K3").Value = 1897.4
$ws.Range("M3").Value = -1783.4

$ws.Range("H64").Value = 1684.3334
$ws.Range("I64").Value = 1621.2
$ws.Range("J64").Value = 2000
$ws.Range("K64").Value = 1621.2
$ws.Range("L64").Value = 2000
$ws.Range("M64").Value = -1396.2
$ws.Range("N64").Value = -2450

$ws.Range("H67").Value = 1684.3334
$ws.Range("I67").Value = 1621.2
$ws.Range("J67").Value = 2000
$ws.Range("K67").Value = 1621.2
$ws.Range("L67").Value = 2000
$ws.Range("M67").Value = -841.2
$ws.Range("N67").Value = -3560

$ws.Range("H86").Value = 3882.8462
$ws.Range("I86").Value = 2746.25
$ws.Range("J86").Value = 5701.4
$ws.Range("K86").Value = 2746.25
$ws.Range("L86").Value = 5701.4
$ws.Range("M86").Value = -1623.25
$ws.Range("N86").Value = -7947.4

$ws.Range("H89").Value = 3882.8462
$ws.Range("I89").Value = 2746.25
$ws.Range("J89").Value = 5701.4
$ws.Range("K89").Value = 13731.25
$ws.Range("L89").Value = 28507
$ws.Range("M89").Value = -8115.25
$ws.Range("N89").Value = -39739

$ws.Range("H107").Value = 890.7778
$ws.Range("I107").Value = 383.375
$ws.Range("J107").Value = 4950
$ws.Range("K107").Value = 383.375
$ws.Range("L107").Value = 4950
$ws.Range("M107").Value = 1536.625
$ws.Range("N107").Value = -8790


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28321.316
$ws.Range("I31").Value = 2747.4375
$ws.Range("K31").Value = 2747.4375
$ws.Range("M31").Value = -2452.4375

$ws.Range("H34").Value = 28321.316
$ws.Range("I34").Value = 2747.4375
$ws.Range("K34").Value = 2747.4375
$ws.Range("M34").Value = -2545.4375

$ws.Range("H58").Value = 3651.087
$ws.Range("I58").Value = 1326.5883
$ws.Range("J58").Value = 10237.167
$ws.Range("K58").Value = 1326.5883
$ws.Range("L58").Value = 10237.167
$ws.Range("M58").Value = -1123.5883
$ws.Range("N58").Value = -10643.167

$ws.Range("H99").Value = 2833.3333
$ws.Range("I99").Value = 2750
$ws.Range("K99").Value = 2750
$ws.Range("M99").Value = -1252

$ws.Range("H105").Value = 7152.4
$ws.Range("I105").Value = 4116
$ws.Range("K105").Value = 4116
$ws.Range("M105").Value = -2369

$ws.Range("H122").Value = 5515.9585
$ws.Range("I122").Value = 1527.1875
$ws.Range("K122").Value = 4581.5625
$ws.Range("M122").Value = -2131.5625

$ws.Range("H126").Value = 2833.3333
$ws.Range("I126").Value = 2750
$ws.Range("K126").Value = 8250
$ws.Range("M126").Value = -5780

$ws.Range("H136").Value = 3651.087
$ws.Range("I136").Value = 1326.5883
$ws.Range("J136").Value = 10237.167
$ws.Range("K136").Value = 3979.7649
$ws.Range("L136").Value = 30711.501
$ws.Range("M136").Value = -1429.7649
$ws.Range("N136").Value = -35811.501


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1630.5333
$ws.Range("I113").Value = 1115.8
$ws.Range("J113").Value = 1887.9
$ws.Range("K113").Value = 3347.4
$ws.Range("L113").Value = 5663.700000000001
$ws.Range("M113").Value = -1177.4
$ws.Range("N113").Value = -10003.7


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 23751.166
$ws.Range("J118").Value = 23751.166
$ws.Range("L118").Value = 23751.166
$ws.Range("N118").Value = -27065.166

$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -16940


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15077.571
$ws.Range("I40").Value = 15077.571
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 15077.571
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -14941.571
$ws.Range("N40").ClearContents()

$ws.Range("H70").Value = 17747.5
$ws.Range("I70").Value = 9247.5
$ws.Range("J70").Value = 21997.5
$ws.Range("K70").Value = 9247.5
$ws.Range("L70").Value = 21997.5
$ws.Range("M70").Value = -8977.5
$ws.Range("N70").Value = -22537.5

$ws.Range("H73").Value = 17747.5
$ws.Range("I73").Value = 9247.5
$ws.Range("J73").Value = 21997.5
$ws.Range("K73").Value = 9247.5
$ws.Range("L73").Value = 21997.5
$ws.Range("M73").Value = -8311.5
$ws.Range("N73").Value = -23869.5

$ws.Range("H82").Value = 6278.357
$ws.Range("J82").Value = 5936.875
$ws.Range("L82").Value = 5936.875
$ws.Range("N82").Value = -6658.875

$ws.Range("H85").Value = 6278.357
$ws.Range("J85").Value = 5936.875
$ws.Range("L85").Value = 5936.875
$ws.Range("N85").Value = -8432.875

$ws.Range("H93").Value = 6274.4
$ws.Range("I93").Value = 5343
$ws.Range("K93").Value = 5343
$ws.Range("M93").Value = -4095


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 25000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 25000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 25000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -26040

$ws.Range("H104").Value = 7500
$ws.Range("J104").Value = 7500
$ws.Range("L104").Value = 7500
$ws.Range("N104").Value = -14488

